$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "No sabe / No responde" value for existing row 27
$ws.Cells.Item(27, 24).Value = 0

# Append the new monthly series row for 01-09-2021
$newRow = 94

# Column A holds a dd-mm-yyyy formatted label that must stay plain text
# (otherwise Excel auto-converts it to a date serial number).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "01-09-2021"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 8345.200000000001
$ws.Cells.Item($newRow, 3).Value = 479.8
$ws.Cells.Item($newRow, 4).Value = 219.9
$ws.Cells.Item($newRow, 5).Value = 830.1
$ws.Cells.Item($newRow, 6).Value = 52.8
$ws.Cells.Item($newRow, 7).Value = 56.7
$ws.Cells.Item($newRow, 8).Value = 756.4
$ws.Cells.Item($newRow, 9).Value = 1607.3
$ws.Cells.Item($newRow, 10).Value = 360.3
$ws.Cells.Item($newRow, 11).Value = 509.7
$ws.Cells.Item($newRow, 12).Value = 206.3
$ws.Cells.Item($newRow, 13).Value = 165.9
$ws.Cells.Item($newRow, 14).Value = 84
$ws.Cells.Item($newRow, 15).Value = 320
$ws.Cells.Item($newRow, 16).Value = 244.1
$ws.Cells.Item($newRow, 17).Value = 507.7
$ws.Cells.Item($newRow, 18).Value = 706.5
$ws.Cells.Item($newRow, 19).Value = 570.2
$ws.Cells.Item($newRow, 20).Value = 77.5
$ws.Cells.Item($newRow, 21).Value = 316.9
$ws.Cells.Item($newRow, 22).Value = 248.2
$ws.Cells.Item($newRow, 23).Value = 0.9
$ws.Cells.Item($newRow, 24).Value = 24
